# Applies a cyclic rotation of species-record data across rows 6, 8 and 9
# of the "Artfynd" sheet:
#   row 9's data -> row 6
#   row 6's data -> row 8
#   row 8's data -> row 9
# Only columns A, B, D, E, F, G, H, Q, R carry data that differs between
# these rows; all other columns (C, I, J, K, P, S, T, U, V, W, Y, Z, AA,
# AB, AD, AE, AF, AG, AT, AW, AX, AY) are identical across the three rows
# and therefore do not need to change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

# Capture the original values of the three affected rows before any
# writes happen, so the rotation doesn't clobber data it still needs.
$orig = @{}
foreach ($r in 6, 8, 9) {
    $orig[$r] = @{}
    foreach ($col in $cols) {
        $orig[$r][$col] = $ws.Range("$col$r").Value2
    }
}

# New row -> source row mapping (cyclic rotation)
$mapping = @{ 6 = 9; 8 = 6; 9 = 8 }

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value2 = $orig[$srcRow][$col]
    }
}
